$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "Proton Afinity" -> "Proton Affinity"
$ws.Range("B1").Value = "Proton Affinity"

# Convert proton affinity values from kcal/mol to Hartree (divide by 627.5),
# applying only to the rows that actually contain data in the diff.
$updates = @{
    "B2"  = -0.591091816
    "B3"  = -0.585828481
    "B4"  = -0.522420553
    "B5"  = -0.56590409
    "B6"  = -0.564128222
    "B7"  = -0.579566449
    "B8"  = -0.577901319
    "B9"  = -0.52503114
    "B10" = -0.560672316
    "B11" = -0.599382834
    "B12" = -0.598052157
    "B13" = -0.582187502
    "B17" = -0.507003379
    "B18" = -0.503959311
    "B19" = -0.575827691
    "B20" = -0.521185772
    "B21" = -0.568031675
    "B31" = -0.560443973
    "B32" = -0.570013997
    "B33" = -0.559223596
    "B34" = -0.572647408
    "B35" = -0.573085917
    "B36" = -0.569191068
    "B37" = -0.576951121
    "B38" = -0.583032453
    "B39" = -0.5327706
    "B40" = -0.582701615
    "B41" = -0.587130592
    "B42" = -0.577526311
    "B43" = -0.576687379
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
